$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure Price column (D) retains its text formatting so values
# like "1.000" or "0.06518" are not re-interpreted as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.040.49"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "1.844.90"
$ws.Range("E3").Value = "  -3.16%  "
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "233.09"
$ws.Range("E5").Value = "  -2.28%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "0.4653"
$ws.Range("E7").Value = "  -2.22%  "
$ws.Range("D8").Value = "0.2801"
$ws.Range("E8").Value = "  -1.05%  "
$ws.Range("D9").Value = "0.06518"
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("D10").Value = "20.15"
$ws.Range("E10").Value = "  +4.41%  "
$ws.Range("D11").Value = "0.07720"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "96.37"
$ws.Range("E12").Value = "  -5.36%  "
$ws.Range("D13").Value = "1.851.88"
$ws.Range("E13").Value = "  -3.22%  "
$ws.Range("E14").Value = "  -2.93%  "
$ws.Range("D15").Value = "0.6622"
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").Value = "282.89"
$ws.Range("E16").Value = "  +4.61%  "
$ws.Range("D17").Value = "30.051.50"
$ws.Range("E17").Value = "  -1.24%  "
$ws.Range("D18").Value = "0.9998"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "12.45"
$ws.Range("E19").Value = "  -1.43%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.096.06"
$ws.Range("E20").Value = "  -2.70%  "
$ws.Range("D21").Value = "5.293"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("E23").Value = "  -3.28%  "
$ws.Range("D24").Value = "6.089"
$ws.Range("E24").Value = "  -2.59%  "
$ws.Range("D25").Value = "166.08"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("D26").Value = "9.210"
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("D27").Value = "18.91"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("D28").Value = "1.909"
$ws.Range("E28").Value = "  -6.90%  "
$ws.Range("D29").Value = "1.354"
$ws.Range("E29").Value = "  -2.28%  "
$ws.Range("D30").Value = "0.09662"
$ws.Range("E30").Value = "  -3.41%  "
$ws.Range("D31").Value = "4.364"
$ws.Range("E31").Value = "  -4.89%  "
$ws.Range("D32").Value = "1.462"
$ws.Range("E32").Value = "  -2.94%  "
$ws.Range("D33").Value = "4.068"
$ws.Range("E33").Value = "  -3.71%  "
$ws.Range("D34").Value = "0.04608"
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("D35").Value = "0.6950"
$ws.Range("E35").Value = "  -4.31%  "
$ws.Range("D36").Value = "1.076"
$ws.Range("E36").Value = "  -3.04%  "
$ws.Range("D37").Value = "0.9991"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "2.704"
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("D39").Value = "0.01840"
$ws.Range("E39").Value = "  -4.05%  "
$ws.Range("D40").Value = "6.282"
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("D41").Value = "2.473"
$ws.Range("E41").Value = "  -5.21%  "
$ws.Range("D42").Value = "71.12"
$ws.Range("E42").Value = "  -3.75%  "
$ws.Range("D43").Value = "0.8530"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").Value = "1.930"
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("D45").Value = "0.9993"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "102.84"
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("D47").Value = "0.4117"
$ws.Range("E47").Value = "  -2.80%  "
$ws.Range("D48").Value = "991.88"
$ws.Range("E48").Value = "  +7.77%  "
$ws.Range("D49").Value = "7.152"
$ws.Range("E49").Value = "  -3.63%  "
$ws.Range("D50").Value = "9.050"
$ws.Range("E50").Value = "  +3.05%  "
$ws.Range("D51").Value = "33.38"
$ws.Range("E51").Value = "  -3.54%  "
